$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2118068965517241
$ws.Range("B3").Value = 0.1805458229957766
$ws.Range("B7").Value = 0.2757242757242757
$ws.Range("B11").Value = 0.216893039049236
$ws.Range("B12").Value = 0.2738805263656158
$ws.Range("B14").Value = 0.2939434534301629
$ws.Range("B15").Value = 0.2832591683289857
$ws.Range("B16").Value = 0.2672127950068266
$ws.Range("B17").Value = 0.2971181376820615
